# "logger and exception added"
#
# The document ends with a short numbered list ("Notes: / - sub bullets").
# Right after the last sub-bullet there is an (empty) paragraph that only
# carries the _GoBack bookmark, followed by a final empty paragraph.
#
# This edit:
#   1. Inserts a new numbered-list paragraph ("Then, write logger file and
#      exception file.") right before the bookmark paragraph.
#   2. Turns the (still empty) bookmark paragraph into a numbered-list item
#      too (numId 1 / ilvl 0), replacing its plain left-indent.
#   3. Inserts a brand-new empty ListParagraph (with the old left-indent)
#      right after the bookmark paragraph, so the trailing blank line is
#      preserved.

$d = $word.ActiveDocument

# --- Locate the paragraph that hosts the "_GoBack" bookmark -----------------
# (COM exposes it even when Bookmarks.Count reports 0 - Word keeps _GoBack
# around internally.) Fall back to "next-to-last paragraph" if that ever
# fails, since that's this document's fixed shape.
$bookmarkParaIndex = -1
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bmStart = $bm.Range.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $pp = $d.Paragraphs.Item($i)
        if ($pp.Range.Start -le $bmStart -and $pp.Range.End -gt $bmStart) {
            $bookmarkParaIndex = $i
        }
    }
} catch {
    $bookmarkParaIndex = -1
}

if ($bookmarkParaIndex -eq -1) {
    $bookmarkParaIndex = $d.Paragraphs.Count - 1
}

$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)

# --- Rebuild the three paragraphs (new bullet / bookmark bullet / blank) ---
# Using InsertXML on the bookmark paragraph's own range lets us set the
# numbering (numId 1, ilvl 0) precisely and drop the old ind="1080" override
# in the same stroke, while preserving the bookmark markup itself.
$xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
       "<w:r><w:t>Then, write logger file and exception file.</w:t></w:r>" +
       "</w:p>" +
       "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
       "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
       "</w:p>" +
       "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
       "<w:pPr><w:pStyle w:val='ListParagraph'/><w:ind w:left='1080'/></w:pPr>" +
       "</w:p>"

$bookmarkPara.Range.InsertXML($xml)
